$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

$ws.Range("B2").Value = "Thu Jan 25 17:29:39 EST 2024"
$ws.Range("B3").Value = "Thu Jan 25 17:29:49 EST 2024"
$ws.Range("B4").Value = "Thu Jan 25 17:29:59 EST 2024"
$ws.Range("B5").Value = "Thu Jan 25 17:30:09 EST 2024"
$ws.Range("B6").Value = "Thu Jan 25 17:30:19 EST 2024"
$ws.Range("B7").Value = "Thu Jan 25 17:30:29 EST 2024"
$ws.Range("B8").Value = "Thu Jan 25 17:30:38 EST 2024"
$ws.Range("B9").Value = "Thu Jan 25 17:30:48 EST 2024"
$ws.Range("B10").Value = "Thu Jan 25 17:30:57 EST 2024"
$ws.Range("B11").Value = "Thu Jan 25 17:31:06 EST 2024"
$ws.Range("B12").Value = "Thu Jan 25 17:31:15 EST 2024"
$ws.Range("B13").Value = "Thu Jan 25 17:31:25 EST 2024"
$ws.Range("B14").Value = "Thu Jan 25 17:31:35 EST 2024"
$ws.Range("B15").Value = "Thu Jan 25 17:31:45 EST 2024"
$ws.Range("B16").Value = "Thu Jan 25 17:31:54 EST 2024"
$ws.Range("B17").Value = "Thu Jan 25 17:32:04 EST 2024"
$ws.Range("B18").Value = "Thu Jan 25 17:32:14 EST 2024"
$ws.Range("B19").Value = "Thu Jan 25 17:32:23 EST 2024"
$ws.Range("B20").Value = "Thu Jan 25 17:32:33 EST 2024"
$ws.Range("B21").Value = "Thu Jan 25 17:32:43 EST 2024"
$ws.Range("B22").Value = "Thu Jan 25 17:32:52 EST 2024"
$ws.Range("B23").Value = "Thu Jan 25 17:33:02 EST 2024"
$ws.Range("B24").Value = "Thu Jan 25 17:33:13 EST 2024"
$ws.Range("B25").Value = "Thu Jan 25 17:33:23 EST 2024"
$ws.Range("B26").Value = "Thu Jan 25 17:33:32 EST 2024"
$ws.Range("B27").Value = "Thu Jan 25 17:33:42 EST 2024"
$ws.Range("B28").Value = "Thu Jan 25 17:33:52 EST 2024"
$ws.Range("B29").Value = "Thu Jan 25 17:34:02 EST 2024"
$ws.Range("B30").Value = "Thu Jan 25 17:34:11 EST 2024"
$ws.Range("B31").Value = "Thu Jan 25 17:34:21 EST 2024"
$ws.Range("B32").Value = "Thu Jan 25 17:34:30 EST 2024"
$ws.Range("B33").Value = "Thu Jan 25 17:34:40 EST 2024"
$ws.Range("B34").Value = "Thu Jan 25 17:34:50 EST 2024"
$ws.Range("B35").Value = "Thu Jan 25 17:34:59 EST 2024"
$ws.Range("B36").Value = "Thu Jan 25 17:35:10 EST 2024"
$ws.Range("B37").Value = "Thu Jan 25 17:35:19 EST 2024"
$ws.Range("B38").Value = "Thu Jan 25 17:35:28 EST 2024"
$ws.Range("B39").Value = "Thu Jan 25 17:35:39 EST 2024"
$ws.Range("B40").Value = "Thu Jan 25 17:35:48 EST 2024"
$ws.Range("B41").Value = "Thu Jan 25 17:35:57 EST 2024"
$ws.Range("B42").Value = "Thu Jan 25 17:36:07 EST 2024"
$ws.Range("B43").Value = "Thu Jan 25 17:36:16 EST 2024"
$ws.Range("B44").Value = "Thu Jan 25 17:36:26 EST 2024"
$ws.Range("B45").Value = "Thu Jan 25 17:36:35 EST 2024"
$ws.Range("B46").Value = "Thu Jan 25 17:36:45 EST 2024"
$ws.Range("B47").Value = "Thu Jan 25 17:36:54 EST 2024"
